# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Chocobo_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1355
$ws.Range("J97").Value = 710
$ws.Range("L97").Value = 2130
$ws.Range("N97").Value = -3122
$ws.Range("H113").Value = 5300.5
$ws.Range("J113").Value = 5300.5
$ws.Range("L113").Value = 5300.5
$ws.Range("N113").Value = -11808.5
$ws.Range("H129").Value = 869.36
$ws.Range("I129").Value = 405.4
$ws.Range("J129").Value = 893.7789299999999
$ws.Range("K129").Value = 1216.2
$ws.Range("L129").Value = 2681.33679
$ws.Range("M129").Value = 3783.8
$ws.Range("N129").Value = -12681.33679
$ws.Range("H132").Value = 42085390
$ws.Range("I132").Value = 50001620
$ws.Range("J132").Value = 2504251.5
$ws.Range("K132").Value = 150004860
$ws.Range("L132").Value = 7512754.5
$ws.Range("M132").Value = -150002330
$ws.Range("N132").Value = -7517814.5
$ws.Range("H137").Value = 1324985
$ws.Range("I137").Value = 1985124.2
$ws.Range("K137").Value = 5955372.6
$ws.Range("M137").Value = -5952822.6
$ws.Range("H141").Value = 119165.766
$ws.Range("I141").Value = 143958.42
$ws.Range("J141").Value = 3466.6667
$ws.Range("K141").Value = 431875.26
$ws.Range("L141").Value = 10400.0001
$ws.Range("M141").Value = -426695.26
$ws.Range("N141").Value = -20760.0001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4434
$ws.Range("I45").Value = 4912
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 4912
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -4535
$ws.Range("N45").Value = -3754
$ws.Range("H74").Value = 4678.852
$ws.Range("I74").Value = 6303.4
$ws.Range("K74").Value = 6303.4
$ws.Range("M74").Value = -5429.4
$ws.Range("H77").Value = 4678.852
$ws.Range("I77").Value = 6303.4
$ws.Range("K77").Value = 31517
$ws.Range("M77").Value = -27149
$ws.Range("H122").Value = 5478.5
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 6971.3335
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 20914.0005
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -25814.0005
$ws.Range("H132").Value = 1957.4839
$ws.Range("I132").Value = 1048.5
$ws.Range("K132").Value = 3145.5
$ws.Range("M132").Value = -615.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 801.5217
$ws.Range("J64").Value = 808.5
$ws.Range("L64").Value = 808.5
$ws.Range("N64").Value = -1258.5
$ws.Range("H67").Value = 801.5217
$ws.Range("J67").Value = 808.5
$ws.Range("L67").Value = 808.5
$ws.Range("N67").Value = -2368.5
$ws.Range("H134").Value = 2584.4688
$ws.Range("I134").Value = 1991.7916
$ws.Range("K134").Value = 5975.3748
$ws.Range("M134").Value = -3440.3748

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2901.8572
$ws.Range("I58").Value = 1728.6227
$ws.Range("K58").Value = 1728.6227
$ws.Range("M58").Value = -1525.6227
$ws.Range("H99").Value = 11115150
$ws.Range("I99").Value = 15386515
$ws.Range("K99").Value = 15386515
$ws.Range("M99").Value = -15385017
$ws.Range("H122").Value = 4944.5713
$ws.Range("I122").Value = 6000
$ws.Range("J122").Value = 4522.4
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 13567.2
$ws.Range("M122").Value = -15550
$ws.Range("N122").Value = -18467.2
$ws.Range("H126").Value = 11115150
$ws.Range("I126").Value = 15386515
$ws.Range("K126").Value = 46159545
$ws.Range("M126").Value = -46157075
$ws.Range("H132").Value = 3735.1155
$ws.Range("I132").Value = 3431.8
$ws.Range("K132").Value = 10295.4
$ws.Range("M132").Value = -7765.400000000001
$ws.Range("H134").Value = 2290.4119
$ws.Range("J134").Value = 4999.8
$ws.Range("L134").Value = 14999.4
$ws.Range("N134").Value = -20069.4
$ws.Range("H136").Value = 2901.8572
$ws.Range("I136").Value = 1728.6227
$ws.Range("K136").Value = 5185.8681
$ws.Range("M136").Value = -2635.8681

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1433.1177
$ws.Range("I113").Value = 1347.5
$ws.Range("J113").Value = 1509.2222
$ws.Range("K113").Value = 1347.5
$ws.Range("L113").Value = 1509.2222
$ws.Range("M113").Value = 822.5
$ws.Range("N113").Value = -5849.2222
$ws.Range("H122").Value = 4838.5
$ws.Range("I122").Value = 2701
$ws.Range("J122").Value = 7688.5
$ws.Range("K122").Value = 8103
$ws.Range("L122").Value = 23065.5
$ws.Range("M122").Value = -5653
$ws.Range("N122").Value = -27965.5
$ws.Range("H126").Value = 3154.01
$ws.Range("I126").Value = 2862.1626
$ws.Range("J126").Value = 4321.4
$ws.Range("K126").Value = 8586.487800000001
$ws.Range("L126").Value = 12964.2
$ws.Range("M126").Value = -6116.487800000001
$ws.Range("N126").Value = -17904.2

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5262.273
$ws.Range("I7").Value = 3500.75
$ws.Range("J7").Value = 6268.857
$ws.Range("K7").Value = 3500.75
$ws.Range("L7").Value = 6268.857
$ws.Range("M7").Value = -3388.75
$ws.Range("N7").Value = -6492.857
$ws.Range("H40").Value = 7930.8887
$ws.Range("I40").Value = 6450
$ws.Range("J40").Value = 8354
$ws.Range("K40").Value = 6450
$ws.Range("L40").Value = 8354
$ws.Range("M40").Value = -6314
$ws.Range("N40").Value = -8626
$ws.Range("H122").Value = 4849.1177
$ws.Range("I122").Value = 2475
$ws.Range("J122").Value = 5579.615
$ws.Range("K122").Value = 7425
$ws.Range("L122").Value = 16738.845
$ws.Range("M122").Value = -4975
$ws.Range("N122").Value = -21638.845
$ws.Range("H126").Value = 5262.273
$ws.Range("I126").Value = 3500.75
$ws.Range("J126").Value = 6268.857
$ws.Range("K126").Value = 10502.25
$ws.Range("L126").Value = 18806.571
$ws.Range("M126").Value = -8032.25
$ws.Range("N126").Value = -23746.571

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4696.8
$ws.Range("I122").Value = 3334
$ws.Range("J122").Value = 5605.3335
$ws.Range("K122").Value = 10002
$ws.Range("L122").Value = 16816.0005
$ws.Range("M122").Value = -7552
$ws.Range("N122").Value = -21716.0005
$ws.Range("H126").Value = 2361.3333
$ws.Range("I126").Value = 1347.5
$ws.Range("K126").Value = 4042.5
$ws.Range("M126").Value = -1572.5
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 8774568
$ws.Range("J132").Value = 33338204
$ws.Range("L132").Value = 100014612
$ws.Range("N132").Value = -100019672
$ws.Range("H136").Value = 4979.05
$ws.Range("I136").Value = 2145.4
$ws.Range("K136").Value = 6436.200000000001
$ws.Range("M136").Value = -3886.200000000001

# WVR row 131: LeveProfitHQ (N131) no longer applicable now that HQ/NQ prices are equal -> clear cell
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N131").ClearContents()
